# Refresh the crypto price/volume table to the latest scraped values.
# (GitHub Actions update, see commit message.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-style updates for the "Price" column (D): these values often look
# numeric (e.g. "1.001"), but the sheet stores them as text, so we force the
# cell to text format, write the value, then restore the original style so
# no visible formatting changes are introduced.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $savedStyle
}

# Row 2
Set-TextValue 'D2' '22.394.97'
$ws.Range('E2').Value = '  -4.09%  '
# Row 3
Set-TextValue 'D3' '1.573.56'
$ws.Range('E3').Value = '  -3.42%  '
# Row 4
Set-TextValue 'D4' '1.001'
$ws.Range('E4').Value = '  -0.07%  '
# Row 5
Set-TextValue 'D5' '1.001'
$ws.Range('E5').Value = '  -0.19%  '
# Row 6
Set-TextValue 'D6' '289.79'
$ws.Range('E6').Value = '  -2.75%  '
# Row 7
Set-TextValue 'D7' '0.3677'
$ws.Range('E7').Value = '  -2.34%  '
# Row 8
$ws.Range('E8').Value = '  -1.23%  '
# Row 9
Set-TextValue 'D9' '0.3388'
$ws.Range('E9').Value = '  -3.48%  '
# Row 10
$ws.Range('E10').Value = '  -2.58%  '
# Row 11
Set-TextValue 'D11' '0.07622'
$ws.Range('E11').Value = '  -4.93%  '
# Row 12
Set-TextValue 'D12' '1.001'
$ws.Range('E12').Value = '  -0.05%  '
# Row 13
Set-TextValue 'D13' '21.35'
$ws.Range('E13').Value = '  -2.16%  '
# Row 14
$ws.Range('E14').Value = '  -3.32%  '
# Row 15
Set-TextValue 'D15' '6.934'
# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.00001138'
$ws.Range('E16').Value = '  -4.34%  '
# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D17' '1.564.38'
$ws.Range('E17').Value = '  -3.96%  '
# Row 18
Set-TextValue 'D18' '89.18'
$ws.Range('E18').Value = '  -6.22%  '
# Row 19
Set-TextValue 'D19' '0.06747'
$ws.Range('E19').Value = '  -2.58%  '
# Row 20
Set-TextValue 'D20' '1.001'
$ws.Range('E20').Value = '  -0.12%  '
# Row 21
Set-TextValue 'D21' '6.251'
$ws.Range('E21').Value = '  -5.71%  '
# Row 22
Set-TextValue 'D22' '16.59'
$ws.Range('E22').Value = '  -3.90%  '
# Row 23
Set-TextValue 'D23' '0.5312'
$ws.Range('E23').Value = '  -6.85%  '
# Row 24
$ws.Range('E24').Value = '  -1.92%  '
# Row 25
Set-TextValue 'D25' '22.404.29'
$ws.Range('E25').Value = '  -4.14%  '
# Row 26
Set-TextValue 'D26' '2.383'
$ws.Range('E26').Value = '  -2.15%  '
# Row 27
Set-TextValue 'D27' '2.979'
$ws.Range('E27').Value = '  +2.23%  '
# Row 28
Set-TextValue 'D28' '20.01'
$ws.Range('E28').Value = '  -3.26%  '
# Row 29
Set-TextValue 'D29' '145.79'
$ws.Range('E29').Value = '  -3.83%  '
# Row 30
Set-TextValue 'D30' '4.980'
$ws.Range('E30').Value = '  -3.73%  '
# Row 31
Set-TextValue 'D31' '125.89'
$ws.Range('E31').Value = '  -4.17%  '
# Row 32
Set-TextValue 'D32' '1.746.55'
$ws.Range('E32').Value = '  -3.66%  '
# Row 33
Set-TextValue 'D33' '1.047'
$ws.Range('E33').Value = '  +8.50%  '
# Row 34
Set-TextValue 'D34' '6.301'
$ws.Range('E34').Value = '  -6.80%  '
# Row 35
Set-TextValue 'D35' '2.002'
$ws.Range('E35').Value = '  -5.41%  '
# Row 36
Set-TextValue 'D36' '10.35'
$ws.Range('E36').Value = '  -7.57%  '
# Row 37
Set-TextValue 'D37' '0.08451'
$ws.Range('E37').Value = '  -2.75%  '
# Row 38
Set-TextValue 'D38' '0.02543'
$ws.Range('E38').Value = '  -5.13%  '
# Row 39
$ws.Range('E39').Value = '  -3.43%  '
# Row 40
Set-TextValue 'D40' '0.06589'
$ws.Range('E40').Value = '  -2.38%  '
# Row 41
Set-TextValue 'D41' '5.567'
$ws.Range('E41').Value = '  -4.29%  '
# Row 42
Set-TextValue 'D42' '11.79'
$ws.Range('E42').Value = '  -7.52%  '
# Row 43
Set-TextValue 'D43' '1.248'
$ws.Range('E43').Value = '  -3.25%  '
# Row 44
Set-TextValue 'D44' '0.6385'
$ws.Range('E44').Value = '  -5.97%  '
# Row 45
Set-TextValue 'D45' '14.30'
$ws.Range('E45').Value = '  -6.93%  '
# Row 46
Set-TextValue 'D46' '0.9997'
$ws.Range('E46').Value = '  -0.17%  '
# Row 47
Set-TextValue 'D47' '0.5997'
$ws.Range('E47').Value = '  -4.64%  '
# Row 48
Set-TextValue 'D48' '3.747'
$ws.Range('E48').Value = '  -3.64%  '
# Row 49
Set-TextValue 'D49' '2.129'
$ws.Range('E49').Value = '  -4.29%  '
# Row 50
Set-TextValue 'D50' '1.261'
$ws.Range('E50').Value = '  +5.53%  '
# Row 51
Set-TextValue 'D51' '123.29'
$ws.Range('E51').Value = '  -2.29%  '
